$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('E2').Value = 'This review article discusses the anesthetic management of patients undergoing laparoscopic surgery, focusing on the physiological effects of pneumoperitoneum, airway management, ventilation strategies, muscle relaxation, and postoperative care. It highlights the importance of understanding the cardiovascular, respiratory, and neurological impacts of laparoscopy to optimize patient safety and outcomes. The article also addresses contraindications, complications, and future directions in laparoscopic anesthesia.'
$ws.Range('F2').Value = 'Review Article'
$ws.Range('G2').Value = 'The article discusses monitoring gas concentrations (CO2) during laparoscopic surgery to manage hypercarbia. It does not explicitly evaluate a specific device but mentions the importance of capnia monitoring.'
$ws.Range('H2').Value = 'The review is based on 51 articles published between 1992 and 2024, focusing on adult populations.'
$ws.Range('I2').Value = 'D3'
$ws.Range('J2').Value = 'A3'
$ws.Range('K2').Value = 'P1'
$ws.Range('L2').Value = 'R3'
$ws.Range("M2").Value = 10
$ws.Range('N2').Value = 'T1'
$ws.Range("S2").Value = 5
$ws.Range('W2').Value = 'exclude'

# Row 3
$ws.Range('E3').Value = 'This review article discusses the recent advancements in medical gas sensing through the integration of artificial intelligence (AI) and nanotechnology. It highlights how AI enhances the accuracy, safety, and efficiency of medical gas monitoring by improving data processing and enabling real-time diagnostics. The article also explores the use of nanostructured materials like metal oxides and carbon-based nanomaterials to increase sensor sensitivity and selectivity, supporting continuous patient monitoring and predictive diagnostics.'
$ws.Range('F3').Value = 'Review'
$ws.Range('G3').Value = 'AI-enabled medical gas sensors using metal oxides and carbon-based nanomaterials'
$ws.Range('H3').Value = 'Not applicable (review article)'
$ws.Range('I3').Value = 'D3'
$ws.Range('J3').Value = 'A2'
$ws.Range('K3').Value = 'P1'
$ws.Range('L3').Value = 'R3'
$ws.Range("M3").Value = 9
$ws.Range('N3').Value = 'T1'
$ws.Range("S3").Value = 7
$ws.Range('W3').Value = 'exclude'

# Row 4
$ws.Range('E4').Value = 'This computer-based simulation study compares avatar-based patient monitoring (Philips Visual Patient Avatar) to conventional monitoring at different viewing distances (8m and 16m). The study found that the avatar-based monitoring significantly improved the perception and recognition of vital signs compared to conventional monitoring, especially with distant vision. The correct recognition rate using the Visual Patient Avatar compared to conventional monitoring at 8 m was increased by 74% and by 51% at 16-meter viewing distance.'
$ws.Range('F4').Value = 'Prospective, single-center, computer-based simulation study'
$ws.Range('G4').Value = 'Philips Visual Patient Avatar'
$ws.Range('H4').Value = '28 anesthesia providers'
$ws.Range('I4').Value = 'D3'
$ws.Range('J4').Value = 'A3'
$ws.Range('K4').Value = 'P1'
$ws.Range('L4').Value = 'R3'
$ws.Range("M4").Value = 10
$ws.Range('W4').Value = 'exclude'

# Row 5
$ws.Range('E5').Value = 'This article introduces the concept of User Interface (UI) Profiles for medical devices within open networked operating rooms, aiming to standardize and enhance the safety and usability of device interactions. It addresses the challenges of interoperability under the IEEE 11073 SDC standards, where devices from different manufacturers communicate, and highlights the need for a standardized method to share HMI requirements. The authors propose a UI Profile that includes design, risk, and process-related UI requirements, which can be used by SDC Consumers during design, usability engineering, and risk management. The article also presents an architectural overview incorporating ISO IEEE 11073-10700 standard requirements and the results of a questionnaire from medical device manufacturers, indicating their views on the necessity, usefulness, and feasibility of UI Profiles.'
$ws.Range('F5').Value = 'Survey'
$ws.Range('G5').Value = 'The article discusses medical devices in general within the context of open networked operating rooms and the application of User Interface (UI) Profiles to enhance their interoperability and usability. It does not focus on a specific device but rather on a concept applicable to various medical devices.'
$ws.Range('H5').Value = 'Eight medical device manufacturers'
$ws.Range('I5').Value = 'D3'
$ws.Range('J5').Value = 'A3'
$ws.Range('K5').Value = 'P3'
$ws.Range('L5').Value = 'R3'
$ws.Range("M5").Value = 12
$ws.Range('N5').Value = 'T1'
$ws.Range("S5").Value = 6
$ws.Range('W5').Value = 'exclude'

# Row 6
$ws.Range('E6').Value = 'This review article discusses the environmental impact of anesthesia, focusing on anesthetic gases as greenhouse gases and strategies to minimize their exposure and waste. It covers low-flow anesthesia (LFA), automated control systems, and other methods to reduce emissions and promote sustainable practices in anesthesiology.'
$ws.Range('F6').Value = 'Review'
$ws.Range('G6').Value = 'Automated low-flow anesthesia machines (e.g., Zeus®, Aisys®, Flow-i®) and vaporizers (e.g., Aladin cassette vaporizer, AnaConDa)'
$ws.Range('H6').Value = 'Not applicable (review article)'
$ws.Range('I6').Value = 'D3'
$ws.Range('J6').Value = 'A2'
$ws.Range('K6').Value = 'P1'
$ws.Range('L6').Value = 'R3'
$ws.Range("M6").Value = 9
$ws.Range('N6').Value = 'T1'
$ws.Range("S6").Value = 7
$ws.Range('W6').Value = 'exclude'

# Row 7
$ws.Range('E7').Value = 'This meta-analysis evaluates the effects of sevoflurane inhalation anesthesia versus propofol intravenous anesthesia on postoperative cognitive function in cancer patients. The study included 41 studies with a total of 4342 patients. The results showed no significant difference in overall POCD incidence between the two anesthetic agents. However, subgroup analysis indicated that sevoflurane was associated with increased POCD at unspecified postoperative times, but decreased POCD at 1 and 3 days postoperatively. MMSE scores were significantly lower in the sevoflurane group at 1, 3, 6, and 12 hours postoperatively, but no significant differences were found at later time points. The authors conclude that sevoflurane and propofol may have differing short-term, but not long-term, negative impacts on cognitive function in cancer patients.'
$ws.Range('F7').Value = 'Meta-analysis'
$ws.Range('G7').Value = 'The article compares sevoflurane inhalation anesthesia to propofol intravenous anesthesia.'
$ws.Range('H7').Value = '4342 patients (2171 sevoflurane, 2171 propofol)'
$ws.Range('I7').Value = 'D3'
$ws.Range('J7').Value = 'A3'
$ws.Range('K7').Value = 'P1'
$ws.Range('L7').Value = 'R3'
$ws.Range("M7").Value = 10
$ws.Range('N7').Value = 'T1'
$ws.Range("S7").Value = 5
$ws.Range('W7').Value = 'exclude'

# Row 8
$ws.Range('E8').Value = 'This secondary analysis of individual patient data from three randomized clinical trials (n=2492) investigates the association between intra-operative chemical (oxygen exposure) and mechanical power (ventilator settings) and postoperative pulmonary complications (PPCs). The study found that both chemical and mechanical power are independently associated with PPCs, demonstrating an additive rather than synergistic effect. An increase of 1 J.min-1 in chemical power was associated with an 8% higher odds of PPCs, while the same increase in mechanical power raised the odds by 5%.'
$ws.Range('F8').Value = 'Retrospective study'
$ws.Range('G8').Value = 'Ventilator (indirectly through mechanical power settings and oxygen administration)'
$ws.Range('H8').NumberFormat = '@'
$ws.Range('H8').Value = '2492'
$ws.Range('I8').Value = 'D3'
$ws.Range('J8').Value = 'A3'
$ws.Range('K8').Value = 'P1'
$ws.Range('L8').Value = 'R3'
$ws.Range("M8").Value = 10
$ws.Range('W8').Value = 'exclude'

# Row 9
$ws.Range('E9').Value = 'This review article discusses the anesthetic management of patients undergoing laparoscopic surgery, focusing on the physiological effects of pneumoperitoneum, airway management, ventilation strategies, muscle relaxation, and postoperative care. It highlights the complications of laparoscopy, including cardiovascular, respiratory, and neurological effects, and provides recommendations for optimizing patient safety and outcomes.'
$ws.Range('F9').Value = 'Review'
$ws.Range('G9').Value = 'The article discusses anesthesia monitoring equipment and techniques used during laparoscopic surgery. While it does not focus on a single device, it mentions capnography for monitoring CO2 levels. The article also mentions transcranial Doppler (TCD) for diagnosis in patients susceptible to ICP increases.'
$ws.Range('H9').Value = 'The review is based on 51 articles and reference texts.'
$ws.Range('I9').Value = 'D3'
$ws.Range('J9').Value = 'A3'
$ws.Range('K9').Value = 'P1'
$ws.Range('L9').Value = 'R3'
$ws.Range("M9").Value = 10
$ws.Range('N9').Value = 'T1'
$ws.Range("S9").Value = 5
$ws.Range('W9').Value = 'exclude'

# Row 10
$ws.Range('E10').Value = 'This review article discusses the recent advancements in medical gas sensing through the integration of artificial intelligence (AI) and nanotechnology. It highlights the importance of accurate medical gas sensing for patient safety, the improvements in sensor technology using nanomaterials like metal oxides and carbon-based nanomaterials, and the role of AI in enhancing data processing, pattern recognition, and real-time monitoring. The review also explores the challenges and limitations of AI-enhanced medical gas sensing and suggests future directions for research and development.'
$ws.Range('F10').Value = 'Review'
$ws.Range('G10').Value = 'AI-enabled medical gas sensors utilizing metal oxides and carbon-based nanomaterials'
$ws.Range('H10').Value = 'Not applicable (review article)'
$ws.Range('I10').Value = 'D3'
$ws.Range('J10').Value = 'A2'
$ws.Range('K10').Value = 'P1'
$ws.Range('L10').Value = 'R3'
$ws.Range("M10").Value = 9
$ws.Range('N10').Value = 'T1'
$ws.Range("S10").Value = 7
$ws.Range('W10').Value = 'exclude'

# Row 11
$ws.Range('E11').Value = 'This prospective, single-center, computer-based simulation study evaluates whether avatar-based patient monitoring (Philips Visual Patient Avatar) improves remote vital sign recognition compared to conventional monitoring at 8 and 16-meter viewing distances. The study found that the avatar-based monitoring significantly improved the perception of vital signs at both distances.'
$ws.Range('F11').Value = 'Prospective, single-center, computer-based simulation study'
$ws.Range('G11').Value = 'Philips Visual Patient Avatar'
$ws.Range('H11').Value = '28 anesthesia providers participated in 112 simulations'
$ws.Range('I11').Value = 'D3'
$ws.Range('J11').Value = 'A3'
$ws.Range('K11').Value = 'P1'
$ws.Range('L11').Value = 'R3'
$ws.Range("M11").Value = 10
$ws.Range('N11').Value = 'T2'
$ws.Range('O11').Value = 'O1'
$ws.Range('P11').Value = 'F1'
$ws.Range('Q11').Value = 'S1'
$ws.Range('R11').Value = 'C1'
$ws.Range("S11").Value = 6
$ws.Range('T11').Value = 'NA'
$ws.Range('U11').Value = 'NA'
$ws.Range('V11').Value = 'NA'
$ws.Range('W11').Value = 'exclude'

# Row 12
$ws.Range('E12').Value = 'This research article discusses the integration of machine-readable user interface (UI) requirements into open networked operating rooms using the IEEE 11073 SDC standards. It proposes extending these standards with UI Profiles provided by medical device manufacturers to enhance usability, patient safety, and operational efficiency. The study includes an architectural overview and a questionnaire to evaluate the feasibility and benefits of UI Profiles, highlighting the need for standardized HMI specifications and addressing the challenges of risk management in open networked solutions.'
$ws.Range('F12').Value = 'Original Research'
$ws.Range('G12').Value = 'The article discusses medical devices that are part of an open networked operating room, focusing on their user interfaces and interoperability based on the IEEE 11073 SDC standards. It does not focus on a specific medical device, but rather on a system-level approach to device integration and usability.'
$ws.Range('H12').Value = 'Eight medical device manufacturers participated in a questionnaire.'
$ws.Range('I12').Value = 'D3'
$ws.Range('J12').Value = 'A3'
$ws.Range('K12').Value = 'P1'
$ws.Range('L12').Value = 'R3'
$ws.Range("M12").Value = 10
$ws.Range('N12').Value = 'T1'
$ws.Range("S12").Value = 6
$ws.Range('W12').Value = 'exclude'

# Row 13
$ws.Range('E13').Value = 'This review article discusses the environmental impact of anesthetic gases, particularly their contribution to global warming. It advocates for the adoption of low-flow anesthesia (LFA) techniques, the use of alternative anesthetic agents, and the implementation of advanced technologies to minimize waste and reduce the carbon footprint of anesthesia practices. The article also touches upon the occupational risks associated with exposure to inhaled anesthetics and emphasizes the importance of training and education for healthcare professionals.'
$ws.Range('F13').Value = 'Review'
$ws.Range('G13').Value = 'Automated Low-Flow Anesthesia Machines (e.g., Zeus®, Aisys®, Flow-i®), Aladin cassette vaporizer, injection vaporizers, Anesthetic Converting Device (AnaConDa)'
$ws.Range('H13').Value = 'N/A'
$ws.Range('I13').Value = 'D3'
$ws.Range('J13').Value = 'A2'
$ws.Range('K13').Value = 'P1'
$ws.Range('L13').Value = 'R3'
$ws.Range("M13").Value = 9
$ws.Range('N13').Value = 'T1'
$ws.Range("S13").Value = 7
$ws.Range('W13').Value = 'exclude'

# Row 14
$ws.Range('E14').Value = 'This meta-analysis evaluates the effects of sevoflurane inhalation anesthesia versus propofol intravenous anesthesia on postoperative cognitive function in cancer patients. The analysis of 41 studies (40 RCTs, 1 cohort) with 4342 patients showed no significant overall difference in POCD incidence between the two anesthetic agents. Subgroup analysis indicated sevoflurane was associated with increased POCD at unspecified postoperative times but decreased POCD at 1 and 3 days postoperatively. MMSE scores were significantly lower in the sevoflurane group at 1, 3, 6, and 12 hours postoperatively, but no significant differences were found at later time points. The authors conclude that sevoflurane and propofol may have differing short-term, but not long-term, impacts on cognitive function in cancer patients.'
$ws.Range('F14').Value = 'Meta-analysis'
$ws.Range('G14').Value = 'The study compares sevoflurane and propofol as anesthetic agents, but does not evaluate the Spacelabs Multigas Module 92518.'
$ws.Range('H14').Value = '4342 patients (2171 sevoflurane, 2171 propofol)'
$ws.Range('I14').Value = 'D3'
$ws.Range('J14').Value = 'A3'
$ws.Range('K14').Value = 'P1'
$ws.Range('L14').Value = 'R3'
$ws.Range("M14").Value = 10
$ws.Range('N14').Value = 'T1'
$ws.Range("S14").Value = 5
$ws.Range('W14').Value = 'exclude'
